# Add "2022-Q4" data to the workbook:
#  - insert a new quarter sheet ("2022-Q4") before the "2022-Q3" sheet, cloning
#    the "2022-Q3" sheet's layout/formatting and filling in the new quarter's
#    fund-holding data.
#  - update the "总计" (summary) sheet: push the existing quarter rows down by
#    one and put the new quarter's totals in the first data row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by copying "2022-Q3" (keeps styles,
#    column widths, page margins, etc. identical to the other quarter sheets)
#    and inserting it immediately before "2022-Q3".
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)

$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# The fund code and numeric-looking figures are stored as TEXT in this table
# (matching the other quarter sheets), so force a text number-format before
# assigning them - otherwise a bare numeric-looking string like "009010"
# gets auto-coerced into the number 9010 (losing the leading zero).
$q4.Range("B2:G2").NumberFormat = "@"
$q4.Range("B2").Value = "009010"
$q4.Range("C2").Value = "华夏兴阳一年持有期混合"
$q4.Range("D2").Value = "27.18"
$q4.Range("E2").Value = "90.65"
$q4.Range("F2").Value = "3.34"
$q4.Range("G2").Value = "0.9078"
$q4.Range("H2").Value = 4

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row for "2022-Q1" (copying
#    the formatting of the existing last data row) and shift the quarter
#    labels/values down so the newest quarter (2022-Q4) is on top.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A4:D4").Copy($total.Range("A5:D5"))
$total.Range("A5").Value = 3
$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.87

$total.Range("B4").Value = "2022-Q2"
$total.Range("D4").Value = 1.17

$total.Range("B3").Value = "2022-Q3"
$total.Range("D3").Value = 0.01

$total.Range("B2").Value = "2022-Q4"
$total.Range("D2").Value = 0.91

# ---------------------------------------------------------------------------
# Restore the original active sheet/selection so the workbook opens on "总计"
# (matches the unchanged bookViews/activeTab in the source file).
# ---------------------------------------------------------------------------
$total.Activate()
$total.Range("A1").Select()
